# Update the "dSF" (column F) values for specific rows to reflect
# repulled data / recalculated means.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -4
    9  = 2
    12 = -2
    14 = 5
    18 = 2
    25 = -11
    26 = -3
    36 = -2
    39 = 5
    42 = -3
    43 = -1
    44 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
